$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.289.52"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.842.95"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.09"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6718"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07437"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.91"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07718"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "1.826.79"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.009"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6725"
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.05"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.154"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "29.262.95"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008311"
$ws.Range("E18").Value = "  +1.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.21"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.164"
$ws.Range("E22").Value = "  -3.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.99"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.709"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1405"
$ws.Range("E26").Value = "  -3.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.515"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.160"
$ws.Range("E29").Value = "  -2.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.069"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05298"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.880"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7529"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.681"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "1.323.66"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01807"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9186"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.974"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.08460"
$ws.Range("E42").Value = "  +15.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.008"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.10"
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("D45").Value = "1.979.57"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.777"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.88"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000121"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.138"
$ws.Range("E50").Value = "  -4.07%  "
$ws.Range("E51").Value = "  -0.03%  "
